$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.146.29'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').Value = '1.656.08'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''218.03'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').Value = '''0.5297'
$ws.Range('E6').Value = '  +1.73%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.2615'
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = '''0.06334'
$ws.Range('E9').Value = '  +1.06%  '
$ws.Range('D10').Value = '''20.45'
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('D11').Value = '''0.07800'
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('D12').Value = '''4.517'
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('D13').Value = '1.656.34'
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').Value = '1.884.30'
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('D15').Value = '''0.5497'
$ws.Range('E15').Value = '  +1.45%  '
$ws.Range('D16').Value = '0.0₅8212'
$ws.Range('E16').Value = '  +1.54%  '
$ws.Range('D17').Value = '''65.41'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '26.144.86'
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').Value = '''4.607'
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('D21').Value = '''191.32'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = '''144.31'
$ws.Range('E25').Value = '  +3.99%  '
$ws.Range('D26').Value = '''0.1229'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('D28').Value = '''15.99'
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('D29').Value = '''1.469'
$ws.Range('E29').Value = '  +4.50%  '
$ws.Range('D30').Value = '''0.05737'
$ws.Range('E30').Value = '  -3.72%  '
$ws.Range('D31').Value = '''1.274'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').Value = '''3.564'
$ws.Range('E32').Value = '  +1.78%  '
$ws.Range('D33').Value = '''3.269'
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('D34').Value = '''1.601'
$ws.Range('E34').Value = '  +3.48%  '
$ws.Range('D35').Value = '''2.803'
$ws.Range('E35').Value = '  +1.87%  '
$ws.Range('D36').Value = '''0.9522'
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('D37').Value = '''2.415'
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('D38').Value = '''0.5744'
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('D39').Value = '''0.01609'
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '''5.799'
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '''0.8507'
$ws.Range('E41').Value = '  +0.98%  '
$ws.Range('D42').Value = '''104.56'
$ws.Range('E42').Value = '  +3.90%  '
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').Value = '1.041.45'
$ws.Range('E44').Value = '  +3.79%  '
$ws.Range('D45').Value = '1.798.37'
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('D46').Value = '''56.92'
$ws.Range('E47').Value = '  +0.56%  '
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('D49').Value = '''7.882'
$ws.Range('E49').Value = '  -0.71%  '
$ws.Range('D50').Value = '''0.05151'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').Value = '''1.443'
$ws.Range('E51').Value = '  -2.62%  '
